# Gantt Chart update: reorder upcoming tasks, rename a couple of tasks,
# mark the completed "Exercise 1" tasks with the "Good" style, and leave
# the workbook focused on the task table (Sheet1) rather than the chart.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Re-order the tasks that sit between "Exercise 2" and "Score Implementation":
# the greeting/progress-tracker/quiz work moves up, a new "Phishing Malware
# Simulation" task is inserted, and Exercise 3 / Exercise 4 move down to
# follow it.
$ws.Range("A15").Value = "Greeting/ Task information "
$ws.Range("A16").Value = "Progress Tracker"
$ws.Range("A17").Value = "Quiz Research and Creation"
$ws.Range("A18").Value = "Quizzes Implementation "
$ws.Range("A19").Value = "Phishing Malware Simulation"
$ws.Range("A20").Value = "Exercise 3 - Safe Web Browsing & Malware Downloads"
$ws.Range("A21").Value = "Excerise 4 - Recognizing Insider Threats"

# "Further Improvements" is renamed to reflect testing feedback, and
# "Potential Further Exercises and Refinement" no longer exists as a task.
$ws.Range("A24").Value = "Improvements from Testing"

# Mark the now-complete "Excerise 1 / Exercise 2" rows with the built-in
# "Good" cell style, matching the earlier completed tasks above them.
$ws.Range("A12:A15").Style = "Good"

# Bring the task table to the front (instead of the Gantt chart sheet)
# and leave the selection on the newly-inserted task row.
$ws.Activate()
$ws.Range("A17").Select()
